$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.269.43'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '1.841.51'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.57'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6277'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07444'
$ws.Range("E8").Value = '  -2.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2893'
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.27'
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07714'
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").Value = '1.842.68'
$ws.Range("E12").Value = '  -2.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.989'
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6766'
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001009'
$ws.Range("E15").Value = '  -4.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.00'
$ws.Range("E16").Value = '  -1.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.143'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").Value = '29.287.35'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.09'
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.27'
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9991'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.374'
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.47'
$ws.Range("E24").Value = '  +0.77%  '
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.407'
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.55'
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06437'
$ws.Range("E28").Value = '  +14.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.390'
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.476'
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.077'
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.044'
$ws.Range("E32").Value = '  -0.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.818'
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.140'
$ws.Range("E34").Value = '  -2.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6946'
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.832'
$ws.Range("E37").Value = '  +3.32%  '
$ws.Range("D38").Value = '1.243.70'
$ws.Range("E38").Value = '  +1.15%  '
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.525'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9099'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9979'
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("D43").Value = '2.003.54'
$ws.Range("E43").Value = '  -14.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.27'
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.15'
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.047'
$ws.Range("E46").Value = '  -2.43%  '
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000116'
$ws.Range("E48").Value = '  -1.44%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.015'
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3941'
$ws.Range("E50").Value = '  -2.13%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.664'
$ws.Range("E51").Value = '  -1.01%  '
